$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.896.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.208.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "81.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.60%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.40"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.01"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.71%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.540.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.220.93"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.805.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.74"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +10.04%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.50"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0872"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.30"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.122"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.51"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.24"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +11.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +19.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.79"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.200"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.18"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0986"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +27.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.441"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.69%  "
